# [Kadastro App] Yeni kayit eklendi: 2944
# Appends the new "2944" record to both the master "Kayitlar" sheet and the
# per-district "Erdemli" sheet (the workbook keeps a duplicate of each
# district's rows on its own tab), mirroring the last existing row's shape.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Kayitlar")
$ws8 = $wb.Worksheets.Item("Erdemli")

$newRecord = @("2944", "2025-09-08", "Erdemli", "1", "ÇAP", "CEMAL TİMUROĞLU (K.Teknisyeni)")

# "Kayitlar" sheet: new data row goes right after the current last row (23 -> 24)
$row1 = 24
for ($col = 1; $col -le 6; $col++) {
    $ws1.Cells.Item($row1, $col).Value = "'" + $newRecord[$col - 1]
}
# Leading apostrophes force text storage (matching every other cell in the
# table, which is text-typed even when it looks numeric/date-like); drop the
# resulting "quote prefix" formatting so the new cells keep the sheet's plain
# default style instead of picking up an extra number format.
$ws1.Range("A${row1}:F${row1}").ClearFormats()

# "Erdemli" sheet: same record, appended after its current last row (22 -> 23)
$row8 = 23
for ($col = 1; $col -le 6; $col++) {
    $ws8.Cells.Item($row8, $col).Value = "'" + $newRecord[$col - 1]
}
$ws8.Range("A${row8}:F${row8}").ClearFormats()
